# Update "想去人数" (want-to-go count) figures on both the "展览" and
# "全部类型" worksheets, which carry duplicate copies of the same table.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 464
    $ws.Range("F3").Value = 3294
    $ws.Range("F4").Value = 85
    $ws.Range("F5").Value = 656
}
